# Updates the cryptos list (Coin, Link, Price, Volume(1h)) to match the
# latest scrape. Mirrors the author's "Updated cryptos list ... with
# GitHub Actions" commit: most rows keep their coin/link but get refreshed
# Price/Volume(1h) figures; rows 8/9 and 31/32 swap which coin occupies
# which rank (OKB now above Cardano, WrappedliquidstakedEther2.0 now above
# Filecoin).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/Volume columns are plain scraped text (e.g. "23.620.42", "0.9989",
# "  +1.78%  "). Excel auto-converts a bare numeric-looking string typed
# into `.Value` into a real Number, which would change the cell's stored
# type from Text. Prefixing numeric-looking values with a leading quote
# (the same trick Excel's own UI uses for "number stored as text") keeps
# them as Text, matching the source data.
function Set-TextValue {
    param($range, [string]$val)
    if ($val -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$') {
        $range.Value = "'" + $val
    } else {
        $range.Value = $val
    }
}

Set-TextValue $ws.Range("B2") 'Bitcoin'
Set-TextValue $ws.Range("C2") 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
Set-TextValue $ws.Range("D2") '23.620.42'
Set-TextValue $ws.Range("E2") '  +1.78%  '
Set-TextValue $ws.Range("B3") 'Ethereum'
Set-TextValue $ws.Range("C3") 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
Set-TextValue $ws.Range("D3") '1.663.05'
Set-TextValue $ws.Range("E3") '  +3.33%  '
Set-TextValue $ws.Range("B4") 'TetherUSD'
Set-TextValue $ws.Range("C4") 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
Set-TextValue $ws.Range("D4") '0.9989'
Set-TextValue $ws.Range("E4") '  -0.33%  '
Set-TextValue $ws.Range("B5") 'USDC'
Set-TextValue $ws.Range("C5") 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextValue $ws.Range("D5") '0.9990'
Set-TextValue $ws.Range("E5") '  -0.30%  '
Set-TextValue $ws.Range("B6") 'BNB'
Set-TextValue $ws.Range("C6") 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextValue $ws.Range("D6") '302.43'
Set-TextValue $ws.Range("E6") '  +0.05%  '
Set-TextValue $ws.Range("B7") 'XRP'
Set-TextValue $ws.Range("C7") 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextValue $ws.Range("D7") '0.3835'
Set-TextValue $ws.Range("E7") '  +1.34%  '
Set-TextValue $ws.Range("B8") 'OKB'
Set-TextValue $ws.Range("C8") 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D8") '51.30'
Set-TextValue $ws.Range("E8") '  -0.88%  '
Set-TextValue $ws.Range("B9") 'Cardano'
Set-TextValue $ws.Range("C9") 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue $ws.Range("D9") '0.3605'
Set-TextValue $ws.Range("E9") '  +2.16%  '
Set-TextValue $ws.Range("B10") 'Polygon'
Set-TextValue $ws.Range("C10") 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range("D10") '1.245'
Set-TextValue $ws.Range("E10") '  +4.02%  '
Set-TextValue $ws.Range("B11") 'Dogecoin'
Set-TextValue $ws.Range("C11") 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws.Range("D11") '0.08192'
Set-TextValue $ws.Range("E11") '  +1.17%  '
Set-TextValue $ws.Range("B12") 'BinanceUSD'
Set-TextValue $ws.Range("C12") 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range("D12") '0.9986'
Set-TextValue $ws.Range("E12") '  -0.36%  '
Set-TextValue $ws.Range("E13") '  +2.45%  '
Set-TextValue $ws.Range("D14") '6.522'
Set-TextValue $ws.Range("E14") '  +2.58%  '
Set-TextValue $ws.Range("D15") '7.520'
Set-TextValue $ws.Range("E15") '  +4.31%  '
Set-TextValue $ws.Range("D16") '0.00001224'
Set-TextValue $ws.Range("E16") '  +1.61%  '
Set-TextValue $ws.Range("D17") '1.660.97'
Set-TextValue $ws.Range("E17") '  +3.28%  '
Set-TextValue $ws.Range("D18") '97.63'
Set-TextValue $ws.Range("E18") '  +3.86%  '
Set-TextValue $ws.Range("D19") '0.06987'
Set-TextValue $ws.Range("E19") '  +1.02%  '
Set-TextValue $ws.Range("D20") '6.846'
Set-TextValue $ws.Range("E20") '  +5.08%  '
Set-TextValue $ws.Range("D21") '17.75'
Set-TextValue $ws.Range("E21") '  +3.60%  '
Set-TextValue $ws.Range("D22") '0.9993'
Set-TextValue $ws.Range("E22") '  -0.24%  '
Set-TextValue $ws.Range("D23") '12.75'
Set-TextValue $ws.Range("E23") '  +3.67%  '
Set-TextValue $ws.Range("D24") '23.635.74'
Set-TextValue $ws.Range("E24") '  +1.88%  '
Set-TextValue $ws.Range("D25") '2.512'
Set-TextValue $ws.Range("E25") '  +0.18%  '
Set-TextValue $ws.Range("E26") '  +0.24%  '
Set-TextValue $ws.Range("D27") '21.27'
Set-TextValue $ws.Range("E27") '  +2.17%  '
Set-TextValue $ws.Range("D28") '152.38'
Set-TextValue $ws.Range("E28") '  +0.91%  '
Set-TextValue $ws.Range("D29") '5.239'
Set-TextValue $ws.Range("E29") '  +0.16%  '
Set-TextValue $ws.Range("D30") '134.18'
Set-TextValue $ws.Range("E30") '  +1.48%  '
Set-TextValue $ws.Range("B31") 'WrappedliquidstakedEther2.0'
Set-TextValue $ws.Range("C31") 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range("D31") '1.841.32'
Set-TextValue $ws.Range("E31") '  +2.84%  '
Set-TextValue $ws.Range("B32") 'Filecoin'
Set-TextValue $ws.Range("C32") 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D32") '7.208'
Set-TextValue $ws.Range("E32") '  +11.02%  '
Set-TextValue $ws.Range("D33") '2.240'
Set-TextValue $ws.Range("E33") '  +6.94%  '
Set-TextValue $ws.Range("D34") '11.94'
Set-TextValue $ws.Range("E34") '  +3.68%  '
Set-TextValue $ws.Range("D35") '1.060'
Set-TextValue $ws.Range("E35") '  -0.58%  '
Set-TextValue $ws.Range("D36") '0.02810'
Set-TextValue $ws.Range("E36") '  +4.04%  '
Set-TextValue $ws.Range("E37") '  +5.95%  '
Set-TextValue $ws.Range("D38") '0.2506'
Set-TextValue $ws.Range("E38") '  +2.40%  '
Set-TextValue $ws.Range("D39") '0.08813'
Set-TextValue $ws.Range("E39") '  +0.91%  '
Set-TextValue $ws.Range("D40") '0.07044'
Set-TextValue $ws.Range("E40") '  +1.52%  '
Set-TextValue $ws.Range("E41") '  +11.74%  '
Set-TextValue $ws.Range("D42") '0.7051'
Set-TextValue $ws.Range("E42") '  +2.89%  '
Set-TextValue $ws.Range("D43") '1.334'
Set-TextValue $ws.Range("E43") '  +1.02%  '
Set-TextValue $ws.Range("D44") '16.12'
Set-TextValue $ws.Range("E44") '  +5.85%  '
Set-TextValue $ws.Range("D45") '0.6568'
Set-TextValue $ws.Range("E45") '  +4.57%  '
Set-TextValue $ws.Range("D46") '2.315'
Set-TextValue $ws.Range("E46") '  +3.48%  '
Set-TextValue $ws.Range("D47") '0.9987'
Set-TextValue $ws.Range("E47") '  -0.26%  '
Set-TextValue $ws.Range("D48") '3.964'
Set-TextValue $ws.Range("E48") '  +0.62%  '
Set-TextValue $ws.Range("D49") '0.07961'
Set-TextValue $ws.Range("E49") '  +1.40%  '
Set-TextValue $ws.Range("D50") '128.22'
Set-TextValue $ws.Range("E50") '  +1.13%  '
Set-TextValue $ws.Range("D51") '1.197'
Set-TextValue $ws.Range("E51") '  +3.08%  '
